$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.415.70"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "2.595.40"
$ws.Range("E3").Value = "  +6.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.35"
$ws.Range("E5").Value = "  +3.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.42"
$ws.Range("E6").Value = "  +2.67%  "

$ws.Range("E7").Value = "  +5.46%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("E9").Value = "  +13.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.40"
$ws.Range("E10").Value = "  +10.13%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +7.23%  "

$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.36"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.18"
$ws.Range("E13").Value = "  +13.60%  "

$ws.Range("D14").Value = "2.990.92"
$ws.Range("E14").Value = "  +6.41%  "

$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").Value = "2.598.55"
$ws.Range("E16").Value = "  +6.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.926"
$ws.Range("E17").Value = "  +8.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.04"
$ws.Range("E18").Value = "  +5.70%  "

$ws.Range("D19").Value = "46.588.11"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("E20").Value = "  +6.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.74"
$ws.Range("E22").Value = "  +7.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.48"
$ws.Range("E23").Value = "  +5.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "274.61"
$ws.Range("E24").Value = "  +10.89%  "

$ws.Range("E25").Value = "  +8.33%  "

$ws.Range("E26").Value = "  +10.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.81"
$ws.Range("E27").Value = "  +39.05%  "

$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.64"
$ws.Range("E30").Value = "  +8.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.31"
$ws.Range("E31").Value = "  +3.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "39.06"
$ws.Range("E32").Value = "  -2.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.36"
$ws.Range("E33").Value = "  +14.00%  "

$ws.Range("E34").Value = "  -7.28%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.83"
$ws.Range("E35").Value = "  +1.64%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0839"
$ws.Range("E36").Value = "  +7.98%  "

$ws.Range("E37").Value = "  +8.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "150.92"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.124"
$ws.Range("E39").Value = "  +8.52%  "

$ws.Range("E40").Value = "  +5.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.28"
$ws.Range("E41").Value = "  +41.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.04"
$ws.Range("E42").Value = "  +4.04%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0331"
$ws.Range("E43").Value = "  +8.61%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.62"
$ws.Range("E44").Value = "  +10.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.07"
$ws.Range("E45").Value = "  +5.10%  "

$ws.Range("D46").Value = "2.131.65"
$ws.Range("E46").Value = "  +7.73%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "93.21"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.50"
$ws.Range("E49").Value = "  +9.84%  "

$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.31"
$ws.Range("E51").Value = "  +7.49%  "

